$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the currency label text on row 7
$ws.Range("B7").Value = "(En Soles)"

# 2. Insert a new column before C for "CENTRO COSTO" (shifts CANTIDAD/TOTAL right)
$ws.Columns("C:C").Insert()

# 3. Insert a new column before E for "P/U" (shifts TOTAL right again)
$ws.Columns("E:E").Insert()

# 4a. Merge the new CENTRO COSTO header cells first (while still blank / unformatted),
#     then copy the exact header formatting (full box border) from the CANTIDAD header so the
#     merge does not get a "split" border treatment.
$ws.Range("C9:C10").Merge()
$ws.Range("D9").Copy()
$ws.Range("C9:C10").PasteSpecial(-4122)

# 4b. For the new P/U column, copy the CANTIDAD header/data box-border formatting onto the
#     still-unmerged cells first, then merge -- this naturally produces the split
#     (top-only / bottom-only) outline that Excel uses for a boxed, merged header pair.
$ws.Range("D9").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E9:E10").Merge()

# 5. Fill in the new header values
$ws.Range("C9").Value = "CENTRO COSTO"
$ws.Range("E9").Value = "P/U"

# 6. Adjust the column widths (closest achievable to target widths)
$ws.Columns("C").ColumnWidth = 18.5
$ws.Columns("D").ColumnWidth = 15.833333333333334
$ws.Columns("E").ColumnWidth = 10.666666666666666

Write-Host "done"
